$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.700.24'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '2.654.36'
$ws.Range('E3').Value = '  +2.53%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '593.51'
$ws.Range('E5').Value = '  +1.94%  '
$ws.Range('D6').Value = '146.24'
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.590'
$ws.Range('E8').Value = '  -1.31%  '
$ws.Range('D9').Value = '0.107'
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('D10').Value = '5.63'
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '27.45'
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('D14').Value = '3.136.70'
$ws.Range('E14').Value = '  +2.68%  '
$ws.Range('D15').Value = '63.732.91'
$ws.Range('E15').Value = '  +1.38%  '
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '2.667.68'
$ws.Range('E17').Value = '  +2.97%  '
$ws.Range('D18').Value = '11.35'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '342.60'
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').Value = '4.35'
$ws.Range('E20').Value = '  -0.92%  '
$ws.Range('D21').Value = '6.77'
$ws.Range('E21').Value = '  +1.18%  '
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').Value = '68.06'
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('E24').Value = '  +12.34%  '
$ws.Range('D25').Value = '1.67'
$ws.Range('E25').Value = '  +5.04%  '
$ws.Range('B26').Value = 'Bittensor'
$ws.Range('C26').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '561.20'
$ws.Range('E26').Value = '  +20.46%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '0.167'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = '8.53'
$ws.Range('E28').Value = '  +2.20%  '
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').Value = '7.96'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('E31').Value = '  +3.12%  '
$ws.Range('D32').Value = '1.79'
$ws.Range('E32').Value = '  +11.60%  '
$ws.Range('E33').Value = '  -1.06%  '
$ws.Range('D34').Value = '175.28'
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = '0.403'
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('D37').Value = '19.09'
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('D38').Value = '4.66'
$ws.Range('E38').Value = '  +2.21%  '
$ws.Range('D39').Value = '1.75'
$ws.Range('E39').Value = '  +2.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '171.80'
$ws.Range('E40').Value = '  +8.18%  '
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D42').Value = '40.46'
$ws.Range('E42').Value = '  +2.66%  '
$ws.Range('E43').Value = '  -0.92%  '
$ws.Range('D44').Value = '21.71'
$ws.Range('E44').Value = '  +3.12%  '
$ws.Range('D45').Value = '0.629'
$ws.Range('E45').Value = '  -1.57%  '
$ws.Range('D46').Value = '0.0548'
$ws.Range('E46').Value = '  +0.41%  '
$ws.Range('D47').Value = '0.0959'
$ws.Range('E47').Value = '  -1.12%  '
$ws.Range('E48').Value = '  +0.53%  '
$ws.Range('D49').Value = '18.63'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('D50').Value = '1.74'
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range('E51').Value = '  -0.83%  '
